# All iOS Store URLs Added - refresh the review counts/ratings table and
# append the newly-scraped apps (Sheet1, A1:D38).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, 'Go Kinetic by Windstream', 4.8, 58675),
    @(1, 'My altafiber', 4.8, 4077),
    @(2, 'Spectrum Access: Enabled Media', 4.6, 120),
    @(3, 'Cox App', 4.6, 377323),
    @(4, 'My Verizon', 4.6, 4700977),
    @(5, 'myCricket App', 4.5, 41715),
    @(6, 'MyDISH Account', 4.6, 314073),
    @(7, 'T-Mobile', 4.8, 2414167),
    @(8, 'Xfinity Mobile', 4.8, 42656),
    @(9, 'Spectrum News: Local Stories', 4.7, 22498),
    @(10, 'My Spectrum', 4.8, 1913188),
    @(11, 'My Sprint Mobile', 4.5, 1047203),
    @(12, 'Verizon My Fios', 4.5, 253099),
    @(13, 'My CenturyLink', 4.4, 112693),
    @(14, 'Visible mobile', 4.4, 66771),
    @(15, 'SpectrumU', 4.3, 161),
    @(16, 'Spectrum TV', 4.7, 393136),
    @(17, 'myAT&amp;T', 4, 3992411),
    @(18, 'Spectrum SportsNet: Live Games', 4.6, 3895),
    @(19, 'MediacomConnect', 3.3, 6380),
    @(20, 'MyMetro', 2.8, 2909),
    @(21, 'MyFrontier', 4.4, 43402),
    @(22, 'Xfinity', 4.4, 832045),
    @(23, 'Google Fiber', 3.5, 201),
    @(24, 'My Viasat', 4.2, 1499),
    @(25, 'Armstrong', 3.9, 13),
    @(26, 'RCN Mobile', 1.7, 25),
    @(27, 'HughesNet Mobile', 1.9, 294),
    @(28, 'HT My Account', 1.4, 17),
    @(29, 'Midco My Account', 1.4, 19),
    @(30, 'Optimum Support', 1.4, 1335),
    @(31, 'UScellular™ – My Account', 4.4, 25203),
    @(32, 'My SECTV', 1.9, 41),
    @(33, 'Optimum TV', 4.6, 12430),
    @(34, 'Breezeline TV', 1.3, 32),
    @(35, 'My Blue Ridge', 4.6, 3081),
    @(36, 'myBuckeye', 2.5, 22)
)

# Rows 23-38 are brand new: clone the formatting of the last
# pre-existing data row (22) down onto each of them first so the new
# "index" cells in column A pick up the same bold/centered/bordered style
# before their values get written below.
for ($r = 23; $r -le 38; $r++) {
    $ws.Range("A22:D22").Copy($ws.Range("A" + $r + ":D" + $r))
}

$rowNum = 2
foreach ($entry in $data) {
    $ws.Cells.Item($rowNum, 1).Value = $entry[0]
    $ws.Cells.Item($rowNum, 2).Value = $entry[1]
    $ws.Cells.Item($rowNum, 3).Value = $entry[2]
    $ws.Cells.Item($rowNum, 4).Value = $entry[3]
    $rowNum = $rowNum + 1
}

